$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2; existing rows 2-15 shift down to 3-16.
$ws.Rows.Item(2).Insert()
# The inserted row inherits the header row's bold/centered formatting;
# strip that so the new row looks like the plain data rows beneath it.
$ws.Range("A2:R2").ClearFormats()

# Populate the newly inserted row 2 with the latest week's data.
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C2").Value = "Metropolitana"
$ws.Range("D2").Value = 44699
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 100112035
$ws.Range("G2").Value = "Bruselas (repollito)"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = 18667
$ws.Range("N2").Value = "$/malla 15 kilos"
$ws.Range("O2").Value = "Provincia de Quillota"
$ws.Range("P2").Value = 1244
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = "Hortaliza"
